# Add 2022-Q1 data.
#
# Strategy: the existing "总计" ("Total") sheet is repurposed to hold the
# newly reported 2022-Q1 fund-holdings detail (this keeps its original
# sheetId/rId), and a brand-new sheet named "总计" is appended at the end
# holding the refreshed totals table (this gets a fresh sheetId/rId),
# matching how the workbook evolves release over release.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Repurpose the current "总计" sheet into the "2022-Q1" detail sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Grab a header cell with the correct style (bold/border/centered) from
# an existing detail sheet so the new headers look the same.
$headerStyleSrc = $wb.Worksheets.Item("2021-Q4").Range("B1")
$headerStyleSrc.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$indexStyleSrc = $wb.Worksheets.Item("2021-Q4").Range("A2")

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q1.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# Columns: B=code  C=name  D=scale  E=stockPosition  F=positionPct  G=marketValue  H=rank
$q1rows = @(
    @("009954","北信瑞丰优选成长股票","0.57","94.37","7.98","0.0455",4),
    @("001829","北信瑞丰中国智造主题灵活配置混合","0.36","94.06","8.57","0.0309",3),
    @("002123","北信瑞丰外延增长主题灵活配置混合","0.17","94.48","8.13","0.0138",3),
    @("004192","招商中证500指数增强A","0.96","94.32","1.26","0.0121",2),
    @("004193","招商中证500指数增强C","0.42","94.32","1.26","0.0053",2),
    @("006195","国金量化多因子股票","0.09","80.71","0.91","0.0008",3)
)

for ($r = 0; $r -lt $q1rows.Length; $r++) {
    $row = 2 + $r
    $rowData = $q1rows[$r]

    $indexStyleSrc.Copy()
    $q1.Cells.Item($row, 1).PasteSpecial(-4122)   # xlPasteFormats
    $q1.Cells.Item($row, 1).Value = $r

    $codeCell = $q1.Cells.Item($row, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $rowData[0]
    $codeCell.Style = "Normal"

    $q1.Cells.Item($row, 3).Value = $rowData[1]

    $scaleCell = $q1.Cells.Item($row, 4)
    $scaleCell.NumberFormat = "@"
    $scaleCell.Value = $rowData[2]
    $scaleCell.Style = "Normal"

    $stockPosCell = $q1.Cells.Item($row, 5)
    $stockPosCell.NumberFormat = "@"
    $stockPosCell.Value = $rowData[3]
    $stockPosCell.Style = "Normal"

    $posPctCell = $q1.Cells.Item($row, 6)
    $posPctCell.NumberFormat = "@"
    $posPctCell.Value = $rowData[4]
    $posPctCell.Style = "Normal"

    $mktValCell = $q1.Cells.Item($row, 7)
    $mktValCell.NumberFormat = "@"
    $mktValCell.Value = $rowData[5]
    $mktValCell.Style = "Normal"

    $q1.Cells.Item($row, 8).Value = $rowData[6]
}

Write-Host "2022-Q1 sheet populated"

# ---------------------------------------------------------------------
# 2) Append a brand-new "总计" sheet with the refreshed totals table.
# ---------------------------------------------------------------------
$totalTmp = $wb.Worksheets.Add()
$totalTmp.Name = "总计"
$lastIndex = $wb.Worksheets.Count
$afterSheet = $wb.Worksheets.Item($lastIndex)
$totalTmp.Move($null, $afterSheet)

# Re-fetch the sheet object by name: the reference obtained before the
# Move() call becomes stale (it keeps pointing at the original sheet
# index instead of following the sheet to its new position).
$total = $wb.Worksheets.Item("总计")

$headerStyleSrc.Copy()
$total.Range("B1:D1").PasteSpecial(-4122)   # xlPasteFormats

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"

# Columns: B=date  C=count  D=marketValue(亿元, numeric)
$totalRows = @(
    @("2022-Q1", 6, 0.11),
    @("2021-Q4", 3, 0.13),
    @("2021-Q3", 6, 0.14),
    @("2021-Q2", 4, 0.09),
    @("2021-Q1", 2, 0.02),
    @("2020-Q4", 1, 0.02)
)

for ($r = 0; $r -lt $totalRows.Length; $r++) {
    $row = 2 + $r
    $rowData = $totalRows[$r]

    $indexStyleSrc.Copy()
    $total.Cells.Item($row, 1).PasteSpecial(-4122)   # xlPasteFormats
    $total.Cells.Item($row, 1).Value = $r

    $total.Cells.Item($row, 2).Value = $rowData[0]
    $total.Cells.Item($row, 3).Value = $rowData[1]
    $total.Cells.Item($row, 4).Value = $rowData[2]
}

Write-Host "总计 sheet populated"
